$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value pairs, derived from the commit diff.
$updates = @(
    @("D2", "37.129.02"),
    @("E2", "  -0.60%  "),
    @("D3", "2.076.15"),
    @("E3", "  -0.91%  "),
    @("E4", "  -0.07%  "),
    @("D5", "253.36"),
    @("E5", "  +0.94%  "),
    @("E6", "  +1.80%  "),
    @("D7", "59.21"),
    @("E7", "  +8.91%  "),
    @("E8", "  -0.01%  "),
    @("D9", "0.392"),
    @("E9", "  +4.47%  "),
    @("D10", "61.49"),
    @("E10", "  -0.45%  "),
    @("E11", "  +7.89%  "),
    @("E12", "  +2.43%  "),
    @("D13", "16.34"),
    @("E13", "  +7.35%  "),
    @("D14", "2.379.58"),
    @("E14", "  -0.91%  "),
    @("D15", "0.819"),
    @("E15", "  -2.24%  "),
    @("E16", "  +6.44%  "),
    @("D17", "2.081.17"),
    @("E17", "  -0.74%  "),
    @("D18", "37.098.41"),
    @("E18", "  -0.53%  "),
    @("D19", "15.84"),
    @("E19", "  +7.77%  "),
    @("D20", "74.86"),
    @("E20", "  +2.74%  "),
    @("D21", "0.0₃0931"),
    @("E21", "  +9.58%  "),
    @("E22", "  +5.07%  "),
    @("D23", "239.50"),
    @("E23", "  -0.71%  "),
    @("E24", "  -0.01%  "),
    @("E25", "  -2.62%  "),
    @("E26", "  +13.85%  "),
    @("D27", "170.42"),
    @("E27", "  -1.14%  "),
    @("D28", "9.36"),
    @("E28", "  +0.94%  "),
    @("D29", "20.40"),
    @("E29", "  -1.38%  "),
    @("E30", "  +2.90%  "),
    @("E31", "  +7.42%  "),
    @("E32", "  +5.96%  "),
    @("E33", "  +3.27%  "),
    @("D34", "4.52"),
    @("E34", "  +8.85%  "),
    @("D35", "0.0910"),
    @("E35", "  +0.14%  "),
    @("E36", "  -0.12%  "),
    @("E37", "  +2.48%  "),
    @("E38", "  +26.15%  "),
    @("D39", "1.78"),
    @("E39", "  -4.23%  "),
    @("E40", "  +2.05%  "),
    @("E41", "  +0.30%  "),
    @("D42", "17.90"),
    @("E42", "  -2.62%  "),
    @("E43", "  +0.16%  "),
    @("D44", "99.29"),
    @("E44", "  +0.30%  "),
    @("B45", "FTXToken"),
    @("C45", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"),
    @("D45", "4.35"),
    @("E45", "  +6.27%  "),
    @("B46", "HuobiToken"),
    @("C46", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"),
    @("D46", "2.84"),
    @("E46", "  +1.42%  "),
    @("E47", "  +13.94%  "),
    @("D48", "2.51"),
    @("E48", "  +7.29%  "),
    @("D49", "1.307.34"),
    @("E49", "  -1.18%  "),
    @("E50", "  -0.10%  "),
    @("D51", "6.94"),
    @("E51", "  -0.78%  ")
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $rng = $ws.Range($cellRef)
    # Force text interpretation so numeric-looking strings (e.g. "253.36")
    # are not coerced into floating point numbers, then drop the explicit
    # number-format style again so the cell keeps its original (default) style.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}
